$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1962.8096
$ws.Range("J17").Value = 2000.95
$ws.Range("L17").Value = 6002.85
$ws.Range("N17").Value = -6338.85

$ws.Range("H42").Value = 1418.7778
$ws.Range("I42").Value = 1804.1428
$ws.Range("J42").Value = 70
$ws.Range("K42").Value = 5412.428400000001
$ws.Range("L42").Value = 210
$ws.Range("M42").Value = -5182.428400000001
$ws.Range("N42").Value = -670

$ws.Range("H55").Value = 244.05882
$ws.Range("I55").Value = 88.9
$ws.Range("J55").Value = 465.7143
$ws.Range("K55").Value = 88.9
$ws.Range("L55").Value = 465.7143
$ws.Range("M55").Value = 125.1
$ws.Range("N55").Value = -893.7143

$ws.Range("H58").Value = 2408.182
$ws.Range("J58").Value = 5063
$ws.Range("L58").Value = 15189
$ws.Range("N58").Value = -15489

$ws.Range("H70").Value = 16250.75
$ws.Range("J70").Value = 16250.75
$ws.Range("L70").Value = 48752.25
$ws.Range("N70").Value = -49292.25

$ws.Range("H73").Value = 16250.75
$ws.Range("J73").Value = 16250.75
$ws.Range("L73").Value = 48752.25
$ws.Range("N73").Value = -50624.25

$ws.Range("H118").Value = 485.63635
$ws.Range("I118").Value = 211.5
$ws.Range("J118").Value = 1216.6666
$ws.Range("K118").Value = 634.5
$ws.Range("L118").Value = 3649.9998
$ws.Range("M118").Value = 1022.5
$ws.Range("N118").Value = -6963.9998

$ws.Range("H129").Value = 887.3333
$ws.Range("I129").Value = 887.3333
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 2661.9999
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 2338.0001
$ws.Range("N129").ClearContents()

$ws.Range("H132").Value = 247835.62
$ws.Range("I132").Value = 281733.03
$ws.Range("K132").Value = 845199.0900000001
$ws.Range("M132").Value = -842669.0900000001

$ws.Range("H137").Value = 15105.857
$ws.Range("I137").Value = 13622.75
$ws.Range("J137").Value = 17083.334
$ws.Range("K137").Value = 40868.25
$ws.Range("L137").Value = 51250.00199999999
$ws.Range("M137").Value = -38318.25
$ws.Range("N137").Value = -56350.00199999999

$ws.Range("H138").Value = 3168.7307
$ws.Range("I138").Value = 1544.0667
$ws.Range("J138").Value = 5384.1816
$ws.Range("K138").Value = 4632.2001
$ws.Range("L138").Value = 16152.5448
$ws.Range("M138").Value = 507.7999
$ws.Range("N138").Value = -26432.5448

$ws.Range("H141").Value = 2018.5625
$ws.Range("I141").Value = 1638.2307
$ws.Range("K141").Value = 4914.6921
$ws.Range("M141").Value = 265.3078999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H27").Value = 2000
$ws.Range("J27").Value = 2000
$ws.Range("L27").Value = 2000
$ws.Range("N27").Value = -2368

$ws.Range("H32").Value = 2026344
$ws.Range("I32").Value = 4698.8726
$ws.Range("K32").Value = 4698.8726
$ws.Range("M32").Value = -4411.8726

$ws.Range("H132").Value = 577969.94
$ws.Range("J132").Value = 86475.82
$ws.Range("L132").Value = 259427.46
$ws.Range("N132").Value = -264487.46

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1652234
$ws.Range("I134").Value = 2399984.5
$ws.Range("K134").Value = 7199953.5
$ws.Range("M134").Value = -7197418.5

$ws.Range("H135").Value = 99995
$ws.Range("J135").Value = 99995
$ws.Range("L135").Value = 99995
$ws.Range("N135").Value = -110135

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 7796.1113
$ws.Range("I62").Value = 9250.833
$ws.Range("K62").Value = 9250.833
$ws.Range("M62").Value = -8626.833

$ws.Range("H65").Value = 7796.1113
$ws.Range("I65").Value = 9250.833
$ws.Range("K65").Value = 46254.165
$ws.Range("M65").Value = -43134.165

$ws.Range("H107").Value = 894.5
$ws.Range("I107").Value = 894.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 894.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1025.5
$ws.Range("N107").ClearContents()

$ws.Range("H134").Value = 55566680
$ws.Range("I134").Value = 83339290
$ws.Range("K134").Value = 250017870
$ws.Range("M134").Value = -250015335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3241.111
$ws.Range("J137").Value = 4164
$ws.Range("L137").Value = 12492
$ws.Range("N137").Value = -22692

$ws.Range("H139").Value = 27780232
$ws.Range("I139").Value = 31252024
$ws.Range("J139").Value = 5899
$ws.Range("K139").Value = 93756072
$ws.Range("L139").Value = 17697
$ws.Range("M139").Value = -93750932
$ws.Range("N139").Value = -27977

$ws.Range("H140").Value = 37502344
$ws.Range("I140").Value = 50001396
$ws.Range("K140").Value = 150004188
$ws.Range("M140").Value = -149999008

$ws.Range("H141").Value = 8019.3335
$ws.Range("I141").Value = 8019.3335
$ws.Range("K141").Value = 24058.0005
$ws.Range("M141").Value = -18878.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 8661.538
$ws.Range("I43").Value = 8661.538
$ws.Range("K43").Value = 8661.538
$ws.Range("M43").Value = -8510.538

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6485.4546
$ws.Range("I122").Value = 6167.625
$ws.Range("K122").Value = 18502.875
$ws.Range("M122").Value = -16052.875

$ws.Range("H132").Value = 2639.5518
$ws.Range("I132").Value = 2366
$ws.Range("J132").Value = 4349.25
$ws.Range("K132").Value = 7098
$ws.Range("L132").Value = 13047.75
$ws.Range("M132").Value = -4568
$ws.Range("N132").Value = -18107.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 20600
$ws.Range("I62").Value = 32750
$ws.Range("K62").Value = 32750
$ws.Range("M62").Value = -32126

$ws.Range("H65").Value = 20600
$ws.Range("I65").Value = 32750
$ws.Range("K65").Value = 163750
$ws.Range("M65").Value = -160630

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H123").Value = 63995
$ws.Range("J123").Value = 63995
$ws.Range("L123").Value = 63995
$ws.Range("N123").Value = -73795

$ws.Range("H126").Value = 3293.6875
$ws.Range("I126").Value = 1919.1538
$ws.Range("K126").Value = 5757.4614
$ws.Range("M126").Value = -3287.4614
